$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 100: 2025-10-20 (serial 45950), 四方坪站 (station "2")
$ws.Range("A100").Value = 45950
$ws.Range("B100").Value = "四方坪站"
$ws.Range("C100").Formula = "=19892/126"
$ws.Range("D100").Formula = "=C100/(24*60)"
$ws.Range("E100").Formula = "=10820.19/126"
$ws.Range("F100").Formula = "=3759.82/126"
$ws.Range("G100").Formula = "=10820.19/(19892/60)"
$ws.Range("H100").Formula = "=427/126"

# New row 101: 2025-10-20 (serial 45950), 高岭站 (station "3")
$ws.Range("A101").Value = 45950
$ws.Range("B101").Value = "高岭站"
$ws.Range("C101").Formula = "=6345/36"
$ws.Range("D101").Formula = "=C101/(24*60)"
$ws.Range("E101").Formula = "=4517.3/36"
$ws.Range("F101").Formula = "=1280.75/36"
$ws.Range("G101").Formula = "=4517.3/(6345/60)"
$ws.Range("H101").Formula = "=157/36"

# Move the selection to match where the author ended up after adding rows
$ws.Range("K99").Select()
